# Seattle Home Pass Types — add a "Match ID" column at the front of the
# sheet (sourced from the SQL schema work mentioned in the commit msg).
#
# Net effect: insert a new column A (pushing the existing A:V data to
# B:W), then populate the new column:
#   - Row 1 (hidden sub-header row): left blank
#   - Row 2 (header row):            "Match ID"
#   - Row 3 (hidden spacer row):     blank
#   - Rows 4-19 (player rows):       34
#   - Row 20 (hidden totals row):    34
# The header/player rows pick up the workbook's existing bold "header"
# font so a brand-new bold-no-border style gets minted, matching the
# rest of the bold cells already on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; this shifts all existing
# data/formatting/merged cells one column to the right automatically.
[void]$ws.Range("A1").EntireColumn.Insert()

# Header cell for the new column.
$ws.Range("A2").Value = "Match ID"
$ws.Range("A2").Font.Bold = $true

# Hidden spacer row underneath the header keeps the same bold styling.
$ws.Range("A3").Font.Bold = $true

# Match ID value (34) for every visible player row.
$ws.Range("A4:A19").Value = 34
$ws.Range("A4:A19").Font.Bold = $true

# Hidden totals row also gets the Match ID value, but keeps default
# (non-bold) formatting.
$ws.Range("A20").Value = 34

# Reflect the selection left behind after filling the new column.
[void]$ws.Range("A2:A19").Select()
